# Atualizacao de bases das ligas, do dia: 17-05-2024 as 13:59
#
# A handful of previously recorded matches had their odds data attached
# to the wrong fixture id. Re-assign the correct data set (columns B and
# E:AB) to each affected row, and refresh a few odds values for two
# still-in-progress fixtures (rows 440-441).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 209 <-> 210 swap ---
$ws.Cells.Item(209, 2).Value = 7404217
$row = 209
$col = 5
$vals = @('Alianza Petrolera', 'Deportivo Pereira', 2, 1, 'H', 1.95, 3.2, 3.75, 1.95, 3.2, 4.75, -0.5, 1.925, 1.875, 2, 1.825, 1.975, 0.95, -1, -1, 0.925, -1, 0.825, -1)
foreach ($v in $vals) {
    $ws.Cells.Item($row, $col).Value = $v
    $col = $col + 1
}

$ws.Cells.Item(210, 2).Value = 7404214
$row = 210
$col = 5
$vals = @('Boyaca Chico', 'Deportivo Cali', 1, 1, 'D', 3.2, 3.1, 2.2, 3.6, 3, 2.25, 0.25, 1.95, 1.9, 2.25, 1.875, 1.975, -1, 2, -1, 0.475, -0.5, -0.5, 0.4875)
foreach ($v in $vals) {
    $ws.Cells.Item($row, $col).Value = $v
    $col = $col + 1
}

# --- Rows 213, 214, 215, 216, 217 cycle ---
$ws.Cells.Item(213, 2).Value = 7404219
$row = 213
$col = 5
$vals = @('Union Magdalena', 'Independiente Medellin', 0, 4, 'A', 3, 3.1, 2.3, 3.6, 3.4, 2.1, 0.25, 2.025, 1.775, 2.5, 1.85, 1.95, -1, -1, 1.1, -1, 0.7749999999999999, 0.8500000000000001, -1)
foreach ($v in $vals) {
    $ws.Cells.Item($row, $col).Value = $v
    $col = $col + 1
}

$ws.Cells.Item(214, 2).Value = 7404215
$row = 214
$col = 5
$vals = @('America de Cali', 'Atletico Bucaramanga', 1, 2, 'A', 1.444, 4.5, 6, 1.363, 5, 7.5, -1.25, 1.775, 2.025, 3, 1.925, 1.875, -1, -1, 6.5, -1, 1.025, 0, 0)
foreach ($v in $vals) {
    $ws.Cells.Item($row, $col).Value = $v
    $col = $col + 1
}

$ws.Cells.Item(215, 2).Value = 7404522
$row = 215
$col = 5
$vals = @('La Equidad', 'Millonarios', 2, 1, 'H', 2.4, 3.1, 2.875, 2.1, 3.1, 3.8, -0.25, 1.75, 2.05, 2, 1.85, 1.95, 1.1, -1, -1, 0.75, -1, 0.8500000000000001, -1)
foreach ($v in $vals) {
    $ws.Cells.Item($row, $col).Value = $v
    $col = $col + 1
}

$ws.Cells.Item(216, 2).Value = 7404260
$row = 216
$col = 5
$vals = @('Atletico Nacional Medellin', 'Deportes Tolima', 2, 3, 'A', 2, 3.25, 3.5, 1.75, 3.6, 4.75, -0.75, 2, 1.8, 2.5, 2, 1.8, -1, -1, 3.75, -1, 0.8, 1, -1)
foreach ($v in $vals) {
    $ws.Cells.Item($row, $col).Value = $v
    $col = $col + 1
}

$ws.Cells.Item(217, 2).Value = 7404213
$row = 217
$col = 5
$vals = @('Jaguares de Cordoba', 'Aguilas Doradas', 0, 1, 'A', 3.25, 3.1, 2.2, 3.6, 3.2, 2.15, 0.25, 1.975, 1.825, 2, 1.75, 2.05, -1, -1, 1.15, -1, 0.825, -1, 1.05)
foreach ($v in $vals) {
    $ws.Cells.Item($row, $col).Value = $v
    $col = $col + 1
}

# --- Row 238 <-> 239 swap ---
$ws.Cells.Item(238, 2).Value = 7528604
$row = 238
$col = 5
$vals = @('Aguilas Doradas', 'Deportivo Cali', 3, 1, 'H', 1.666, 3.75, 5, 1.363, 5, 9, -1.25, 1.825, 1.975, 2.75, 1.9, 1.9, 0.363, -1, -1, 0.825, -1, 0.8999999999999999, -1)
foreach ($v in $vals) {
    $ws.Cells.Item($row, $col).Value = $v
    $col = $col + 1
}

$ws.Cells.Item(239, 2).Value = 7528136
$row = 239
$col = 5
$vals = @('Millonarios', 'Atletico Nacional Medellin', 0, 1, 'A', 1.85, 3.3, 4.5, 1.85, 3.5, 4.2, -0.5, 1.875, 1.975, 2.5, 2.05, 1.8, -1, -1, 3.2, -1, 0.9750000000000001, -1, 0.8)
foreach ($v in $vals) {
    $ws.Cells.Item($row, $col).Value = $v
    $col = $col + 1
}

# --- Rows 425, 426, 427, 428, 429, 430 cycle ---
$ws.Cells.Item(425, 2).Value = 7658987
$row = 425
$col = 5
$vals = @('Deportivo Cali', 'Junior', 0, 0, 'D', 2.7, 3.25, 2.4, 3.2, 3.1, 2.4, 0.25, 1.8, 2.05, 2.25, 1.975, 1.875, -1, 2.1, -1, 0.4, -0.5, -1, 0.875)
foreach ($v in $vals) {
    $ws.Cells.Item($row, $col).Value = $v
    $col = $col + 1
}

$ws.Cells.Item(426, 2).Value = 7658989
$row = 426
$col = 5
$vals = @('Jaguares de Cordoba', 'Independiente Santa Fe', 1, 0, 'H', 3, 3.2, 2.3, 3.4, 3.6, 2.05, 0.25, 2, 1.8, 2.5, 1.8, 2, 2.4, -1, -1, 1, -1, -1, 1)
foreach ($v in $vals) {
    $ws.Cells.Item($row, $col).Value = $v
    $col = $col + 1
}

$ws.Cells.Item(427, 2).Value = 7658990
$row = 427
$col = 5
$vals = @('Millonarios', 'Boyaca Chico', 3, 0, 'H', 1.4, 4.2, 7, 1.4, 4.5, 8.5, -1.25, 1.95, 1.9, 2.5, 1.975, 1.875, 0.3999999999999999, -1, -1, 0.95, -1, 0.9750000000000001, -1)
foreach ($v in $vals) {
    $ws.Cells.Item($row, $col).Value = $v
    $col = $col + 1
}

$ws.Cells.Item(428, 2).Value = 7658988
$row = 428
$col = 5
$vals = @('Envigado FC', 'Independiente Medellin', 0, 1, 'A', 4.2, 3.4, 1.8, 5.25, 3.6, 1.7, 0.75, 1.925, 1.875, 2.25, 1.775, 2.025, -1, -1, 0.7, -0.5, 0.4375, -1, 1.025)
foreach ($v in $vals) {
    $ws.Cells.Item($row, $col).Value = $v
    $col = $col + 1
}

$ws.Cells.Item(429, 2).Value = 7658915
$row = 429
$col = 5
$vals = @('Once Caldas', 'America de Cali', 0, 0, 'D', 2.3, 3, 3.1, 2.3, 3.2, 3.3, -0.25, 1.975, 1.875, 2.25, 2.025, 1.825, -1, 2.2, -1, -0.5, 0.4375, -1, 0.825)
foreach ($v in $vals) {
    $ws.Cells.Item($row, $col).Value = $v
    $col = $col + 1
}

$ws.Cells.Item(430, 2).Value = 7658914
$row = 430
$col = 5
$vals = @('La Equidad', 'Deportivo Pereira', 0, 2, 'A', 2, 3.1, 3.75, 2.25, 3.2, 3.3, -0.25, 1.925, 1.875, 2, 1.825, 1.975, -1, -1, 2.3, -1, 0.875, 0, 0)
foreach ($v in $vals) {
    $ws.Cells.Item($row, $col).Value = $v
    $col = $col + 1
}

# --- Odds refresh for in-progress fixtures (rows 440-441) ---
$ws.Cells.Item(440, 14).Value = 3.1   # N440
$ws.Cells.Item(440, 17).Value = 2.05  # Q440
$ws.Cells.Item(440, 18).Value = 1.8   # R440

$ws.Cells.Item(441, 13).Value = 2.8    # M441
$ws.Cells.Item(441, 15).Value = 2.75   # O441
$ws.Cells.Item(441, 17).Value = 1.925  # Q441
$ws.Cells.Item(441, 18).Value = 1.925  # R441
